# Insert a new weekly data row for "Tomate / Larga vida / Primera / Región del Maule"
# as row 302, pushing the existing rows 302-321 down to 303-322.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("302").Insert()

$ws.Cells.Item(302, 1).Value = 7
$ws.Cells.Item(302, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(302, 3).Value = "Ñuble"
$ws.Cells.Item(302, 4).Value = 44516
$ws.Cells.Item(302, 5).Value = 16
$ws.Cells.Item(302, 6).Value = 100112020
$ws.Cells.Item(302, 7).Value = "Tomate"
$ws.Cells.Item(302, 8).Value = "Larga vida"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 300
$ws.Cells.Item(302, 11).Value = 13000
$ws.Cells.Item(302, 12).Value = 14000
$ws.Cells.Item(302, 13).Value = 13500
$ws.Cells.Item(302, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(302, 15).Value = "Región del Maule"
$ws.Cells.Item(302, 16).Value = 900
$ws.Cells.Item(302, 17).Value = 15
$ws.Cells.Item(302, 18).Value = "Hortaliza"

$ws.Cells.Item(302, 4).NumberFormat = $ws.Cells.Item(303, 4).NumberFormat
